# Stat Sync & Fixed superfluous error messages
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Characters")

# Fixed superfluous error messages:
# - H12 ("Magnetic" passive) had a leftover/erroneous note; clear it.
$ws.Cells.Item(12, 8).Value = $null

# - H2 ("Fire" passive) repeated info already implied elsewhere; trim it down.
$ws.Cells.Item(2, 8).Value = "Nearby burning enemies give strength and hp regeneration"

# Stat Sync: add a baseline stat-sync row (Speed/Health reference values),
# matching the number formatting/alignment already used by the I:L columns.
$ws.Range("I15:J15").Copy()
$ws.Range("I16:J16").PasteSpecial(-4122)
$ws.Cells.Item(16, 9).Value = 50
$ws.Cells.Item(16, 10).Value = 50
